$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '37.116.19'
Set-TextValue 'E2' '  -0.12%  '

Set-TextValue 'D3' '2.049.83'
Set-TextValue 'E3' '  -0.59%  '

Set-TextValue 'E4' '  +0.01%  '

Set-TextValue 'D5' '248.61'
Set-TextValue 'E5' '  -0.70%  '

Set-TextValue 'D6' '0.665'
Set-TextValue 'E6' '  -1.91%  '

Set-TextValue 'D7' '59.31'
Set-TextValue 'E7' '  -0.11%  '

Set-TextValue 'D8' '0.999'
Set-TextValue 'E8' '  -0.05%  '

Set-TextValue 'D9' '0.384'
Set-TextValue 'E9' '  +0.89%  '

Set-TextValue 'D10' '0.0787'
Set-TextValue 'E10' '  -1.66%  '

Set-TextValue 'E11' '  +1.09%  '

Set-TextValue 'D12' '15.80'
Set-TextValue 'E12' '  +3.81%  '

Set-TextValue 'D13' '2.344.47'
Set-TextValue 'E13' '  -0.74%  '

Set-TextValue 'E14' '  +1.70%  '

Set-TextValue 'E15' '  +7.28%  '

Set-TextValue 'D16' '2.047.02'
Set-TextValue 'E16' '  -0.79%  '

Set-TextValue 'D17' '17.88'
Set-TextValue 'E17' '  +23.27%  '

Set-TextValue 'D18' '37.077.80'
Set-TextValue 'E18' '  -0.10%  '

Set-TextValue 'D19' '75.18'
Set-TextValue 'E19' '  -0.22%  '

Set-TextValue 'D20' '0.0₃0900'
Set-TextValue 'E20' '  -2.63%  '

Set-TextValue 'D21' '5.35'
Set-TextValue 'E21' '  -0.65%  '

Set-TextValue 'D22' '237.64'
Set-TextValue 'E22' '  -0.71%  '

Set-TextValue 'E23' '  +0.03%  '

Set-TextValue 'E24' '  +0.45%  '

Set-TextValue 'D25' '169.15'
Set-TextValue 'E25' '  -1.45%  '

Set-TextValue 'E26' '  +7.63%  '

Set-TextValue 'D27' '9.38'
Set-TextValue 'E27' '  +1.85%  '

Set-TextValue 'D28' '20.04'
Set-TextValue 'E28' '  -1.25%  '

Set-TextValue 'E29' '  -0.42%  '

Set-TextValue 'D30' '1.12'
Set-TextValue 'E30' '  +4.13%  '

Set-TextValue 'D31' '4.78'
Set-TextValue 'E31' '  +3.20%  '

Set-TextValue 'D32' '0.0624'
Set-TextValue 'E32' '  -1.55%  '

Set-TextValue 'D33' '4.52'
Set-TextValue 'E33' '  +2.51%  '

Set-TextValue 'D34' '0.0898'
Set-TextValue 'E34' '  +0.98%  '

Set-TextValue 'E35' '  -0.05%  '

Set-TextValue 'E36' '  -2.96%  '

Set-TextValue 'E37' '  -1.18%  '

Set-TextValue 'E38' '  -0.68%  '

Set-TextValue 'E39' '  -4.64%  '

Set-TextValue 'D40' '3.17'
Set-TextValue 'E40' '  +13.38%  '

Set-TextValue 'D41' '5.04'
Set-TextValue 'E41' '  +12.28%  '

Set-TextValue 'B42' 'InjectiveProtocol'
Set-TextValue 'C42' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D42' '17.47'
Set-TextValue 'E42' '  -5.21%  '

Set-TextValue 'B43' 'VeChain'
Set-TextValue 'C43' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D43' '0.0221'
Set-TextValue 'E43' '  -2.10%  '

Set-TextValue 'D44' '1.14'
Set-TextValue 'E44' '  -1.37%  '

Set-TextValue 'D45' '96.18'
Set-TextValue 'E45' '  -1.58%  '

Set-TextValue 'D46' '2.47'
Set-TextValue 'E46' '  -1.77%  '

Set-TextValue 'E47' '  -0.38%  '

Set-TextValue 'D48' '1.283.81'
Set-TextValue 'E48' '  -1.63%  '

Set-TextValue 'D49' '6.77'
Set-TextValue 'E49' '  -1.88%  '

Set-TextValue 'D50' '2.230.40'
Set-TextValue 'E50' '  -0.92%  '

Set-TextValue 'D51' '3.54'
Set-TextValue 'E51' '  -21.46%  '
